$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''27.519.56'
$ws.Range('E2').Value = '  +5.70%  '
$ws.Range('D3').Value = '''1.823.86'
$ws.Range('E3').Value = '  +6.29%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''346.25'
$ws.Range('E5').Value = '  +4.85%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('D7').Value = '''0.3830'
$ws.Range('E7').Value = '  +3.83%  '
$ws.Range('D8').Value = '''0.3534'
$ws.Range('E8').Value = '  +6.32%  '
$ws.Range('D9').Value = '''49.40'
$ws.Range('E9').Value = '  -0.84%  '
$ws.Range('D10').Value = '''1.240'
$ws.Range('E10').Value = '  +5.03%  '
$ws.Range('D11').Value = '''0.07808'
$ws.Range('D12').Value = '''1.002'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '''22.32'
$ws.Range('E13').Value = '  +11.51%  '
$ws.Range('D14').Value = '''6.658'
$ws.Range('E14').Value = '  +6.47%  '
$ws.Range('D15').Value = '''7.287'
$ws.Range('E15').Value = '  +5.49%  '
$ws.Range('D16').Value = '''1.821.04'
$ws.Range('E16').Value = '  +6.31%  '
$ws.Range('D17').Value = '''0.00001131'
$ws.Range('E17').Value = '  +5.14%  '
$ws.Range('D18').Value = '''0.06745'
$ws.Range('E18').Value = '  +1.75%  '
$ws.Range('D19').Value = '''86.63'
$ws.Range('E19').Value = '  +5.73%  '
$ws.Range('D21').Value = '''17.78'
$ws.Range('E21').Value = '  +8.99%  '
$ws.Range('D22').Value = '''6.566'
$ws.Range('E22').Value = '  +8.40%  '
$ws.Range('E23').Value = '  +2.76%  '
$ws.Range('D24').Value = '''27.523.86'
$ws.Range('E24').Value = '  +5.92%  '
$ws.Range('D25').Value = '''2.463'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').Value = '''2.698'
$ws.Range('E26').Value = '  +8.78%  '
$ws.Range('D27').Value = '''22.34'
$ws.Range('E27').Value = '  +16.08%  '
$ws.Range('D28').Value = '''1.516'
$ws.Range('E28').Value = '  +15.60%  '
$ws.Range('D29').Value = '''153.77'
$ws.Range('E29').Value = '  +2.65%  '
$ws.Range('D30').Value = '''2.021.55'
$ws.Range('E30').Value = '  +6.12%  '
$ws.Range('D31').Value = '''136.82'
$ws.Range('E31').Value = '  +6.21%  '
$ws.Range('D32').Value = '''6.413'
$ws.Range('E32').Value = '  +7.73%  '
$ws.Range('D33').Value = '''4.010'
$ws.Range('E33').Value = '  -2.37%  '
$ws.Range('D34').Value = '''14.18'
$ws.Range('E34').Value = '  +9.98%  '
$ws.Range('D35').Value = '''0.08803'
$ws.Range('E35').Value = '  +3.25%  '
$ws.Range('D36').Value = '''1.703'
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').Value = '''5.686'
$ws.Range('E37').Value = '  +6.38%  '
$ws.Range('D38').Value = '''0.7105'
$ws.Range('E38').Value = '  +15.49%  '
$ws.Range('D39').Value = '''0.2295'
$ws.Range('E39').Value = '  +7.82%  '
$ws.Range('D40').Value = '''0.06568'
$ws.Range('E40').Value = '  +6.11%  '
$ws.Range('D41').Value = '''0.02430'
$ws.Range('E41').Value = '  +6.43%  '
$ws.Range('D42').Value = '''9.065'
$ws.Range('E42').Value = '  +6.68%  '
$ws.Range('D43').Value = '''1.305'
$ws.Range('E43').Value = '  +1.90%  '
$ws.Range('D44').Value = '''14.88'
$ws.Range('E44').Value = '  +3.11%  '
$ws.Range('D45').Value = '''0.6645'
$ws.Range('E45').Value = '  +13.39%  '
$ws.Range('E46').Value = '  +0.20%  '
$ws.Range('D47').Value = '''4.050'
$ws.Range('E47').Value = '  +5.65%  '
$ws.Range('D48').Value = '''2.204'
$ws.Range('E48').Value = '  +9.91%  '
$ws.Range('D49').Value = '''133.20'
$ws.Range('E49').Value = '  +4.69%  '
$ws.Range('D50').Value = '''0.07373'
$ws.Range('E50').Value = '  +1.60%  '
$ws.Range('D51').Value = '''81.02'
$ws.Range('E51').Value = '  +5.23%  '
